$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '61.404.51'
Set-TextValue 2 5 '  +0.92%  '

# Row 3
Set-TextValue 3 4 '3.407.92'
Set-TextValue 3 5 '  +3.02%  '

# Row 4
Set-TextValue 4 4 '0.999'
Set-TextValue 4 5 '  -0.47%  '

# Row 5
Set-TextValue 5 4 '406.86'
Set-TextValue 5 5 '  -0.62%  '

# Row 6
Set-TextValue 6 4 '129.17'
Set-TextValue 6 5 '  +15.47%  '

# Row 7
Set-TextValue 7 4 '0.603'
Set-TextValue 7 5 '  +6.93%  '

# Row 8
Set-TextValue 8 4 '3.400.45'
Set-TextValue 8 5 '  +3.19%  '

# Row 9
Set-TextValue 9 5 '  -0.03%  '

# Row 10
Set-TextValue 10 4 '0.674'
Set-TextValue 10 5 '  +9.48%  '

# Row 11
Set-TextValue 11 4 '0.128'
Set-TextValue 11 5 '  +19.50%  '

# Row 12
Set-TextValue 12 4 '42.27'
Set-TextValue 12 5 '  +9.55%  '

# Row 13
Set-TextValue 13 5 '  -0.70%  '

# Row 14
Set-TextValue 14 4 '3.931.26'
Set-TextValue 14 5 '  +1.85%  '

# Row 15
Set-TextValue 15 4 '8.57'
Set-TextValue 15 5 '  +5.03%  '

# Row 16
Set-TextValue 16 4 '19.70'
Set-TextValue 16 5 '  +3.66%  '

# Row 17
Set-TextValue 17 4 '3.393.98'
Set-TextValue 17 5 '  +0.13%  '

# Row 18
Set-TextValue 18 4 '11.74'
Set-TextValue 18 5 '  +12.91%  '

# Row 19
Set-TextValue 19 4 '61.262.87'
Set-TextValue 19 5 '  +0.66%  '

# Row 20
Set-TextValue 20 4 '1.02'
Set-TextValue 20 5 '  +3.06%  '

# Row 21
Set-TextValue 21 4 '0.0000134'
Set-TextValue 21 5 '  +21.34%  '

# Row 22
Set-TextValue 22 4 '3.25'
Set-TextValue 22 5 '  +0.86%  '

# Row 23
Set-TextValue 23 4 '82.63'
Set-TextValue 23 5 '  +12.99%  '

# Row 24
Set-TextValue 24 4 '13.07'
Set-TextValue 24 5 '  +7.46%  '

# Row 25
Set-TextValue 25 4 '308.21'
Set-TextValue 25 5 '  +4.98%  '

# Row 26
Set-TextValue 26 4 '3.23'
Set-TextValue 26 5 '  +5.64%  '

# Row 27
Set-TextValue 27 4 '8.61'
Set-TextValue 27 5 '  +16.23%  '

# Row 28
Set-TextValue 28 2 'LEO'
Set-TextValue 28 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 28 4 '4.66'
Set-TextValue 28 5 '  +4.18%  '

# Row 29
Set-TextValue 29 2 'EthereumClassic'
Set-TextValue 29 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 29 4 '29.76'
Set-TextValue 29 5 '  +4.78%  '

# Row 30
Set-TextValue 30 4 '7.51'
Set-TextValue 30 5 '  +1.68%  '

# Row 31
Set-TextValue 31 4 '0.174'
Set-TextValue 31 5 '  +3.76%  '

# Row 32
Set-TextValue 32 4 '0.116'
Set-TextValue 32 5 '  +6.19%  '

# Row 33
Set-TextValue 33 4 '11.72'
Set-TextValue 33 5 '  +5.25%  '

# Row 34
Set-TextValue 34 4 '42.76'
Set-TextValue 34 5 '  +6.65%  '

# Row 35
Set-TextValue 35 4 '2.58'
Set-TextValue 35 5 '  +7.85%  '

# Row 36
Set-TextValue 36 5 '  +0.43%  '

# Row 37
Set-TextValue 37 4 '0.0486'
Set-TextValue 37 5 '  +3.04%  '

# Row 38
Set-TextValue 38 4 '52.21'
Set-TextValue 38 5 '  +0.48%  '

# Row 39
Set-TextValue 39 4 '0.994'
Set-TextValue 39 5 '  -0.76%  '

# Row 40
Set-TextValue 40 4 '3.43'
Set-TextValue 40 5 '  +4.11%  '

# Row 41
Set-TextValue 41 4 '3.01'
Set-TextValue 41 5 '  -2.56%  '

# Row 42
Set-TextValue 42 4 '0.125'
Set-TextValue 42 5 '  +5.76%  '

# Row 43
Set-TextValue 43 4 '1.98'
Set-TextValue 43 5 '  +5.94%  '

# Row 44
Set-TextValue 44 4 '135.63'
Set-TextValue 44 5 '  -1.01%  '

# Row 45
Set-TextValue 45 2 'NEARProtocol'
Set-TextValue 45 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 45 4 '3.95'
Set-TextValue 45 5 '  +5.18%  '

# Row 46
Set-TextValue 46 2 'TheGraph'
Set-TextValue 46 3 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 46 4 '0.285'
Set-TextValue 46 5 '  +3.41%  '

# Row 47
Set-TextValue 47 2 'Celestia'
Set-TextValue 47 3 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 47 4 '17.01'
Set-TextValue 47 5 '  +5.44%  '

# Row 48
Set-TextValue 48 5 '  +2.77%  '

# Row 49
Set-TextValue 49 4 '21.90'
Set-TextValue 49 5 '  +5.66%  '

# Row 50
Set-TextValue 50 2 'Maker'
Set-TextValue 50 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 50 4 '2.148.75'
Set-TextValue 50 5 '  +1.00%  '

# Row 51
Set-TextValue 51 2 'RocketPoolETH'
Set-TextValue 51 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 51 4 '3.740.93'
Set-TextValue 51 5 '  +0.29%  '
